# Regenerate s_vals data to filter save games.
# Updates columns B:G for rows 2-10 with new computed values.
# Column G is the row-wise sum of B+C+D+E (F/"Win" is not part of the sum).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(1.455362044514542, 0.306821227259698, 3.537761648806719, 0.4942365360607697, 0, 5.794181456641729)
    3  = @(0.1190320826869504, 0.04071648406533734, 3.537761648806719, 0.4942365360607697, 0, 4.191746751619776)
    4  = @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 0, 3.754798637575387)
    5  = @(0.6606524410359556, 10.34677158129881, 0.7527432677738641, 10.19245300693656, 0, 21.95262029704519)
    6  = @(0.6606524410359556, 1.655778082260271, 3.537761648806719, 10.19245300693656, 1, 16.0466451790395)
    7  = @(0.6606524410359556, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 0, 1.611132179096228)
    8  = @(0.01293466051926884, 0.04071648406533734, 0.7527432677738641, 0.4942365360607697, 1, 1.30063094841924)
    9  = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 1, 6.189590430959694)
    10 = @(3.286832544864788, 0.306821227259698, 3.537761648806719, 10.19245300693656, 0, 17.32386842786776)
}

$cols = @("B", "C", "D", "E", "F", "G")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i]
    }
}
